$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D3: remove the "+ Requests erklären (get…)" suffix, keep trailing double space
$ws.Range("D3").Value = "Screencast zur Datenaquise (Wo sind gute Geospatial Web Services zu finden)  "

# Update D4: "und" -> "des", move the "+ Requests erklären (get…)" part to the end without trailing space
$ws.Range("D4").Value = "Vorstellung des Funktionsrahmen der verschiedenen Geospatial Web Services darstellen (WMS, WFS, etc.) + Requests erklären (get…)"

# Set width for new column E
$ws.Columns.Item(5).ColumnWidth = 76

# Apply the same wrap-text formatting used by the other data columns
$ws.Range("E2:E7").WrapText = $true

# Fill E column values (only E3 gets text, rest stay empty)
$ws.Range("E3").Value = "Ist das notwendig? Oder lieber einfach einige Links dazu geben oder in Video 1 kur mit aufnehmen "

# Adjust row heights
$ws.Rows.Item(3).RowHeight = 29
$ws.Rows.Item(4).RowHeight = 29
$ws.Rows.Item(6).RowHeight = 16.5

# Update the selected cell to mirror the author's final selection
$ws.Range("D12").Select() | Out-Null
